# Update the public EPEX Spot prices workbook with the latest day's data.

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value (e.g. an ISO date string) into a cell
# without Excel's automatic "looks like a date" conversion turning it into
# a date serial number.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Sheet "Prix Spot": add column G (20-jun) -----------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("F1").Copy()
$wsSpot.Range("G1").PasteSpecial(-4122)
$wsSpot.Range("G1").Value = "20-jun"

$gValues = @{
    2  = 55.73
    3  = 21.73
    4  = 25.89
    5  = 25.45
    6  = 13.47
    7  = 30.38
    8  = 37.45
    9  = 57.51
    10 = 67.17
    11 = 38.01
    12 = 4.48
    13 = 0
    14 = -0.01
    15 = -0.01
    16 = -0.01
    17 = 0.05
    18 = 12.8
    19 = 50.91
    20 = 94.78
    21 = 115
    22 = 126.68
    23 = 115.67
    24 = 122.86
    25 = 114.9
}

foreach ($row in $gValues.Keys) {
    $wsSpot.Cells.Item($row, 7).Value = $gValues[$row]
}

# --- Sheet "Gaz": add row 4 (2025-06-18) ----------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

Set-TextValue $wsGaz.Cells.Item(4, 1) "2025-06-18"
$wsGaz.Cells.Item(4, 2).Value = 38.45
$wsGaz.Cells.Item(4, 3).Value = 10800
$wsGaz.Cells.Item(4, 4).Value = 37.806

# --- Sheet "CO2": add row 4 (2025-06-18) ----------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

Set-TextValue $wsCo2.Cells.Item(4, 1) "2025-06-18"
$wsCo2.Cells.Item(4, 2).Value = 73.45
